$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.955.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.644.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5226'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2605'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06365'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07679'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.645.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.869.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5509'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8230'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.977.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.698'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '188.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.256'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1249'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.393'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.394'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05898'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.263'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.398'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.643'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9910'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.393'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.750'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5633'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.46%  '
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.871'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8522'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.036.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.71'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.792.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.55%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.004'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.041'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05146'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.76%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4217'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.881'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.10%  '
